$wb = $excel.ActiveWorkbook

# --- "sets" sheet: set 3's home_points changes from 2 to 4 ---
$sets = $wb.Worksheets.Item("sets")
$sets.Cells.Item(4, 4).Value = 4

# --- "rallies" sheet: append two new rally rows (66 and 67) ---
$rallies = $wb.Worksheets.Item("rallies")

# Row 66 -> rally_id 65
$rallies.Cells.Item(66, 1).Value = 65
$rallies.Cells.Item(66, 2).Value = 1
$rallies.Cells.Item(66, 3).Value = 3
$rallies.Cells.Item(66, 4).Value = 3
$rallies.Cells.Item(66, 5).Value = "NOS"
$rallies.Cells.Item(66, 6).Value = ""
$rallies.Cells.Item(66, 7).Value = 3
$rallies.Cells.Item(66, 8).Value = "MEIO"
$rallies.Cells.Item(66, 9).Value = "PONTO"
$rallies.Cells.Item(66, 10).Value = "NOS"
$rallies.Cells.Item(66, 11).Value = 3
$rallies.Cells.Item(66, 12).Value = 0
$rallies.Cells.Item(66, 13).Value = "1 3 m"
$rallies.Cells.Item(66, 14).Value = "FRENTE"
$rallies.Cells.Item(66, 15).Value = "FRENTE"
$rallies.Cells.Item(66, 16).Value = "FRENTE"

# Row 67 -> rally_id 66
$rallies.Cells.Item(67, 1).Value = 66
$rallies.Cells.Item(67, 2).Value = 1
$rallies.Cells.Item(67, 3).Value = 3
$rallies.Cells.Item(67, 4).Value = 4
$rallies.Cells.Item(67, 5).Value = "NOS"
$rallies.Cells.Item(67, 6).Value = ""
$rallies.Cells.Item(67, 7).Value = 3
$rallies.Cells.Item(67, 8).Value = "MEIO"
$rallies.Cells.Item(67, 9).Value = "PONTO"
$rallies.Cells.Item(67, 10).Value = "NOS"
$rallies.Cells.Item(67, 11).Value = 4
$rallies.Cells.Item(67, 12).Value = 0
$rallies.Cells.Item(67, 13).Value = "1 3 m"
$rallies.Cells.Item(67, 14).Value = "FRENTE"
$rallies.Cells.Item(67, 15).Value = "FRENTE"
$rallies.Cells.Item(67, 16).Value = "FRENTE"
